$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 102, pushing the existing rows 102-124 down to 103-125
# (weekly price data: a new week's record is inserted at the top of this block).
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with this week's record.
$ws.Cells.Item(102, 1).Value = 7
$ws.Cells.Item(102, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(102, 3).Value = "Ñuble"
$ws.Cells.Item(102, 4).Value = 44476
$ws.Cells.Item(102, 5).Value = 16
$ws.Cells.Item(102, 6).Value = 100112017
$ws.Cells.Item(102, 7).Value = "Apio"
$ws.Cells.Item(102, 8).Value = "Americana (o)"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 300
$ws.Cells.Item(102, 11).Value = 7500
$ws.Cells.Item(102, 12).Value = 8000
$ws.Cells.Item(102, 13).Value = 7750
$ws.Cells.Item(102, 14).Value = "`$/docena de matas"
$ws.Cells.Item(102, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(102, 16).Value = 1292
$ws.Cells.Item(102, 17).Value = 6
$ws.Cells.Item(102, 18).Value = "Hortaliza"
